$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.722.27"
$ws.Range("E2").Value = "  +0.53%  "
$ws.Range("D3").Value = "'1.859.24"
$ws.Range("E3").Value = "  +0.44%  "
$ws.Range("D4").Value = "'1.019"
$ws.Range("E4").Value = "  -1.07%  "
$ws.Range("D5").Value = "'320.93"
$ws.Range("E5").Value = "  -0.05%  "
$ws.Range("E6").Value = "  -0.97%  "
$ws.Range("D7").Value = "'0.4370"
$ws.Range("E7").Value = "  -0.46%  "
$ws.Range("E8").Value = "  -0.26%  "
$ws.Range("D9").Value = "'0.07417"
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("D10").Value = "'0.8840"
$ws.Range("E10").Value = "  +0.74%  "
$ws.Range("D11").Value = "'21.58"
$ws.Range("E11").Value = "  -0.35%  "
$ws.Range("D12").Value = "'1.858.57"
$ws.Range("E12").Value = "  -0.49%  "
$ws.Range("E13").Value = "  +0.94%  "
$ws.Range("D14").Value = "'5.492"
$ws.Range("E14").Value = "  -0.46%  "
$ws.Range("D15").Value = "'0.07141"
$ws.Range("E15").Value = "  -0.83%  "
$ws.Range("D16").Value = "'87.89"
$ws.Range("E16").Value = "  +5.84%  "
$ws.Range("E17").Value = "  -1.07%  "
$ws.Range("D18").Value = "'0.000009037"
$ws.Range("E18").Value = "  -0.20%  "
$ws.Range("D19").Value = "'1.018"
$ws.Range("E19").Value = "  -0.99%  "
$ws.Range("D20").Value = "'15.42"
$ws.Range("E20").Value = "  -0.09%  "
$ws.Range("D21").Value = "'27.721.68"
$ws.Range("E21").Value = "  +0.42%  "
$ws.Range("D22").Value = "'5.284"
$ws.Range("E22").Value = "  +0.28%  "
$ws.Range("E23").Value = "  -1.96%  "
$ws.Range("D24").Value = "'2.103.77"
$ws.Range("E24").Value = "  +1.25%  "
$ws.Range("D25").Value = "'2.035"
$ws.Range("E25").Value = "  +6.14%  "
$ws.Range("D26").Value = "'157.19"
$ws.Range("E26").Value = "  -0.42%  "
$ws.Range("E27").Value = "  -0.33%  "
$ws.Range("D28").Value = "'5.431"
$ws.Range("E28").Value = "  +2.76%  "
$ws.Range("D29").Value = "'1.992"
$ws.Range("E29").Value = "  +0.89%  "
$ws.Range("D30").Value = "'121.17"
$ws.Range("E30").Value = "  +3.29%  "
$ws.Range("D31").Value = "'0.09053"
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("D32").Value = "'1.218"
$ws.Range("E32").Value = "  +1.29%  "
$ws.Range("D33").Value = "'0.7707"
$ws.Range("E33").Value = "  +1.12%  "
$ws.Range("D34").Value = "'3.037"
$ws.Range("E34").Value = "  +5.35%  "
$ws.Range("D35").Value = "'4.564"
$ws.Range("E35").Value = "  +0.78%  "
$ws.Range("E36").Value = "  -0.91%  "
$ws.Range("D37").Value = "'1.139"
$ws.Range("E37").Value = "  -1.05%  "
$ws.Range("D38").Value = "'0.01982"
$ws.Range("E38").Value = "  +0.23%  "
$ws.Range("D39").Value = "'0.05308"
$ws.Range("E39").Value = "  -0.07%  "
$ws.Range("D40").Value = "'2.874"
$ws.Range("E40").Value = "  +1.90%  "
$ws.Range("D41").Value = "'0.5180"
$ws.Range("E41").Value = "  +0.38%  "
$ws.Range("D42").Value = "'6.964"
$ws.Range("E42").Value = "  +2.66%  "
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("D44").Value = "'8.711"
$ws.Range("E44").Value = "  +2.29%  "
$ws.Range("E45").Value = "  +1.89%  "
$ws.Range("D46").Value = "'110.28"
$ws.Range("E46").Value = "  +1.26%  "
$ws.Range("D47").Value = "'1.715"
$ws.Range("E47").Value = "  +0.15%  "
$ws.Range("D48").Value = "'0.4733"
$ws.Range("E48").Value = "  +1.69%  "
$ws.Range("D49").Value = "'1.020"
$ws.Range("E49").Value = "  -0.95%  "
$ws.Range("D50").Value = "'0.06480"
$ws.Range("E50").Value = "  +1.29%  "
$ws.Range("D51").Value = "'1.848"
$ws.Range("E51").Value = "  -0.17%  "
